$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Remove the second sheet (Hoja2) - its single cell of data is no longer needed.
$excel.DisplayAlerts = $false
[void]$wb.Worksheets.Item("Hoja2").Delete()

# Replace the old NombreProducto/Precio/Marca table with the new
# NombreProducto/Genero product list.
$ws.Range("A1").Value = "NombreProducto"
$ws.Range("B1").Value = "Genero"

$ws.Range("A2").Value = "Reebok Hiit Tr Dynred"
$ws.Range("B2").Value = "Hombre"
$ws.Range("B2").NumberFormat = "@"

$ws.Range("A3").Value = "Advanced Trainette"
$ws.Range("A3").Font.Name = "Calibri"
$ws.Range("B3").Value = "Mujer"

$ws.Range("A4").Value = "Downshifter 9"
$ws.Range("B4").Value = "Hombre"

$ws.Range("A5").Value = "Lebron Witness 4"
$ws.Range("B5").Value = "Hombre"
$ws.Range("B5").Font.Name = "Calibri"

$ws.Range("A6").Value = "Air Max 720"
$ws.Range("B6").Value = "Mujer"
$ws.Range("B6").Font.Name = "Calibri"

$ws.Range("A7").Value = "Court Borough Low 2"
$ws.Range("B7").Value = "Niños"
$ws.Range("B7").Font.Name = "Calibri"

# The old sheet had a third column (Marca) that's now gone.
$ws.Range("C1").Value = $null
$ws.Range("C2").Value = $null

# Leave the selection the way it ended up after filling in the table.
[void]$ws.Rows("4:4").Select()
